$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = "'88.315.44"
$ws.Range('E2').Value = "'  -1.74%  "

$ws.Range('D3').Value = "'3.081.26"
$ws.Range('E3').Value = "'  -3.60%  "

$ws.Range('E4').Value = "'  -0.09%  "

$ws.Range('D5').Value = "'210.23"
$ws.Range('E5').Value = "'  -2.74%  "

$ws.Range('D6').Value = "'623.85"
$ws.Range('E6').Value = "'  +0.31%  "

$ws.Range('D7').Value = "'0.373"
$ws.Range('E7').Value = "'  -3.41%  "

$ws.Range('D8').Value = "'0.821"
$ws.Range('E8').Value = "'  +17.73%  "

$ws.Range('E9').Value = "'  -0.01%  "

$ws.Range('D10').Value = "'3.078.30"
$ws.Range('E10').Value = "'  -3.56%  "

$ws.Range('D11').Value = "'0.592"
$ws.Range('E11').Value = "'  +3.87%  "

$ws.Range('E12').Value = "'  +0.44%  "

$ws.Range('D13').Value = "'0.0000239"
$ws.Range('E13').Value = "'  -5.58%  "

$ws.Range('D14').Value = "'5.29"
$ws.Range('E14').Value = "'  -1.49%  "

$ws.Range('D15').Value = "'88.008.17"
$ws.Range('E15').Value = "'  -1.80%  "

$ws.Range('D16').Value = "'3.645.86"
$ws.Range('E16').Value = "'  -3.75%  "

$ws.Range('D17').Value = "'31.67"
$ws.Range('E17').Value = "'  -3.08%  "

$ws.Range('D18').Value = "'3.068.44"
$ws.Range('E18').Value = "'  -4.25%  "

$ws.Range('D19').Value = "'3.21"
$ws.Range('E19').Value = "'  -3.83%  "

$ws.Range('D20').Value = "'0.0000213"
$ws.Range('E20').Value = "'  -8.67%  "

$ws.Range('D21').Value = "'13.11"
$ws.Range('E21').Value = "'  -1.42%  "

$ws.Range('D22').Value = "'421.68"
$ws.Range('E22').Value = "'  -2.24%  "

$ws.Range('E23').Value = "'  -3.92%  "

$ws.Range('D24').Value = "'4.84"
$ws.Range('E24').Value = "'  -3.75%  "

$ws.Range('D25').Value = "'5.46"
$ws.Range('E25').Value = "'  +7.35%  "

$ws.Range('D26').Value = "'81.64"
$ws.Range('E26').Value = "'  +8.72%  "

$ws.Range('D27').Value = "'11.40"
$ws.Range('E27').Value = "'  -1.40%  "

$ws.Range('D28').Value = "'3.235.89"
$ws.Range('E28').Value = "'  -3.60%  "

$ws.Range('D29').Value = "'0.999"
$ws.Range('E29').Value = "'  +0.02%  "

$ws.Range('D30').Value = "'1.08"
$ws.Range('E30').Value = "'  +8.23%  "

$ws.Range('D31').Value = "'0.155"
$ws.Range('E31').Value = "'  -0.43%  "

$ws.Range('D32').Value = "'8.04"
$ws.Range('E32').Value = "'  -4.25%  "

$ws.Range('D33').Value = "'506.28"
$ws.Range('E33').Value = "'  -5.13%  "

$ws.Range('D34').Value = "'3.57"
$ws.Range('E34').Value = "'  -10.85%  "

$ws.Range('D35').Value = "'6.62"
$ws.Range('E35').Value = "'  -2.95%  "

$ws.Range('D36').Value = "'1.81"
$ws.Range('E36').Value = "'  -3.75%  "

$ws.Range('E37').Value = "'  -2.77%  "

$ws.Range('D38').Value = "'22.23"
$ws.Range('E38').Value = "'  +0.28%  "

$ws.Range('B39').Value = 'WhiteBITCoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D39').Value = "'22.22"
$ws.Range('E39').Value = "'  -0.34%  "

$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').Value = "'0.129"
$ws.Range('E40').Value = "'  +2.52%  "

$ws.Range('E41').Value = "'  +0.31%  "

$ws.Range('E42').Value = "'  -0.02%  "

$ws.Range('B43').Value = 'Monero'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D43').Value = "'149.14"
$ws.Range('E43').Value = "'  -0.98%  "

$ws.Range('B44').Value = 'PolygonEcosystemToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D44').Value = "'0.358"
$ws.Range('E44').Value = "'  -3.12%  "

$ws.Range('D45').Value = "'1.82"
$ws.Range('E45').Value = "'  -4.61%  "

$ws.Range('D46').Value = "'0.135"
$ws.Range('E46').Value = "'  +9.39%  "

$ws.Range('D47').Value = "'43.46"
$ws.Range('E47').Value = "'  +0.71%  "

$ws.Range('D48').Value = "'0.0655"
$ws.Range('E48').Value = "'  +9.96%  "

$ws.Range('D49').Value = "'156.23"
$ws.Range('E49').Value = "'  -8.68%  "

$ws.Range('E50').Value = "'  -3.92%  "

$ws.Range('D51').Value = "'1.17"
$ws.Range('E51').Value = "'  -4.82%  "
